# Auto-generated edit script: strip footnote markers "[n]" and flatten
# embedded line breaks to single spaces in vaccine/brand-name cells.
# Also fixes the two "Afluria Quadrivalent" cells on the Adult Influenza
# sheet that previously held a duplicate two-line variant of the string.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Pediatric VFC Vaccine ")
$ws.Range("A2").Value2 = 'DTaP '
$ws.Range("A3").Value2 = 'DTaP '
$ws.Range("A4").Value2 = 'DTaP '
$ws.Range("A5").Value2 = 'DTaP-IPV '
$ws.Range("A6").Value2 = 'DTaP-IPV '
$ws.Range("A7").Value2 = 'DTaP-IPV '
$ws.Range("A8").Value2 = 'DTaP-Hep B-IPV '
$ws.Range("A9").Value2 = 'DTaP-IP-HI '
$ws.Range("A10").Value2 = 'e-IPV '
$ws.Range("A11").Value2 = 'Hepatitis A Pediatric '
$ws.Range("A12").Value2 = 'Hepatitis A Pediatric '
$ws.Range("A13").Value2 = 'Hepatitis A Pediatric '
$ws.Range("A14").Value2 = 'Hepatitis A Pediatric '
$ws.Range("A15").Value2 = 'Hepatitis A-Hepatitis B 18 only '
$ws.Range("A16").Value2 = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range("A17").Value2 = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range("B17").Value2 = 'Recombivax HB'
$ws.Range("A18").Value2 = 'Hepatitis B  Pediatric/Adolescent'
$ws.Range("B18").Value2 = 'Recombivax HB'
$ws.Range("A19").Value2 = 'Hib '
$ws.Range("A20").Value2 = 'Hib '
$ws.Range("A21").Value2 = 'Hib '
$ws.Range("A22").Value2 = 'HPV - Human Papillomavirus 9-valent '
$ws.Range("A23").Value2 = 'MENB - Meningococcal Group B '
$ws.Range("A24").Value2 = 'MENB - Meningococcal Group B '
$ws.Range("A25").Value2 = 'MENB - Meningococcal Group B '
$ws.Range("A26").Value2 = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range("A27").Value2 = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range("A28").Value2 = 'Measles, Mumps and Rubella (MMR) '
$ws.Range("A29").Value2 = 'MMR/Varicella '
$ws.Range("A30").Value2 = 'Pneumococcal 13-valent  (Pediatric)'
$ws.Range("A32").Value2 = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range("A33").Value2 = 'Rotavirus, Live, Oral, Pentavalent '
$ws.Range("A34").Value2 = 'Rotavirus, Live, Oral, Oral '
$ws.Range("A35").Value2 = 'Tetanus and Diphtheria Toxoids '
$ws.Range("A36").Value2 = 'Tetanus and Diphtheria Toxoids '
$ws.Range("A37").Value2 = 'Tetanus and Diphtheria Toxoids '
$ws.Range("A38").Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A39").Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A40").Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A41").Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A42").Value2 = 'Varicella '

$ws = $wb.Worksheets.Item("Adult Vaccine ")
$ws.Range("A2").Value2 = 'Hepatitis A-Adult '
$ws.Range("A3").Value2 = 'Hepatitis A-Adult '
$ws.Range("A4").Value2 = 'Hepatitis A Adult '
$ws.Range("A5").Value2 = 'Hepatitis A-Hepatitis B Adult '
$ws.Range("A6").Value2 = 'Hepatitis B-Adult '
$ws.Range("A7").Value2 = 'Hepatitis B-Adult '
$ws.Range("A8").Value2 = 'HPV-Human Papillomavirus 9 Valent '
$ws.Range("A9").Value2 = 'Measles, Mumps,  Rubella-Adult '
$ws.Range("A10").Value2 = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range("A11").Value2 = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws.Range("A12").Value2 = 'MENB - Meningococcal Group B '
$ws.Range("A13").Value2 = 'MENB - Meningococcal Group B '
$ws.Range("A14").Value2 = 'MENB - Meningococcal Group B '
$ws.Range("A15").Value2 = 'Pneumococcal 13-valent  (Adult)'
$ws.Range("A18").Value2 = 'Tetanus and Diphtheria Toxoids '
$ws.Range("A19").Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A20").Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A21").Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A22").Value2 = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$ws.Range("A23").Value2 = 'Varicella-Adult '

$ws = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$ws.Range("A2").Value2 = 'Influenza  (Age 6 months and older)'
$ws.Range("B2").Value2 = 'Fluzone Quadrivalent'
$ws.Range("A3").Value2 = 'Influenza  (Age 6-35 months)'
$ws.Range("B3").Value2 = 'Fluzone Quadrivalent Pediatric dose'
$ws.Range("A4").Value2 = 'Influenza  (Age 36 months and older)'
$ws.Range("B4").Value2 = 'Fluzone Quadrivalent'
$ws.Range("A5").Value2 = 'Influenza  (Age 36 months and older)'
$ws.Range("B5").Value2 = 'Fluzone Quadrivalent'
$ws.Range("A6").Value2 = 'Influenza  (Age 6 months and older)'
$ws.Range("B6").Value2 = 'Fluarix Quadrivalent'
$ws.Range("A7").Value2 = 'Influenza  (Age 6 months and older)'
$ws.Range("B7").Value2 = 'FluLaval Quadrivalent'
$ws.Range("A8").Value2 = 'Influenza  (Age 6 months and older)'
$ws.Range("B8").Value2 = 'FluLaval Quadrivalent'
$ws.Range("A9").Value2 = 'Influenza  (Age 4 years and older)'
$ws.Range("A10").Value2 = 'Influenza  (Age 4 years and older)'
$ws.Range("A11").Value2 = 'Influenza  (Age 5 years and older)'
$ws.Range("A12").Value2 = 'Influenza  (Age 5 years and older)'
$ws.Range("A13").Value2 = 'Influenza  Live, Intranasal (Age 2-49 years)'
$ws.Range("B13").Value2 = 'FluMist Quadrivalent'

$ws = $wb.Worksheets.Item("Adult Influenza Vaccine ")
$ws.Range("A2").Value2 = 'Influenza  (Age 6 months and older)'
$ws.Range("B2").Value2 = 'Fluzone Quadrivalent'
$ws.Range("A3").Value2 = 'Influenza  (Age 36 months and older)'
$ws.Range("B3").Value2 = 'Fluzone Quadrivalent'
$ws.Range("A4").Value2 = 'Influenza  (Age 36 months and older)'
$ws.Range("B4").Value2 = 'Fluzone Quadrivalent'
$ws.Range("A5").Value2 = 'Influenza  (Age 6 months and older)'
$ws.Range("B5").Value2 = 'Fluarix Quadrivalent'
$ws.Range("A6").Value2 = 'Influenza  (Age 6 months and older)'
$ws.Range("B6").Value2 = 'FluLaval Quadrivalent'
$ws.Range("A7").Value2 = 'Influenza  (Age 6 months and older)'
$ws.Range("B7").Value2 = 'FluLaval Quadrivalent'
$ws.Range("A8").Value2 = 'Influenza  (Age 4 years and older)'
$ws.Range("A9").Value2 = 'Influenza  (Age 4 years and older)'
$ws.Range("A10").Value2 = 'Influenza  (Age 5 years and older)'
$ws.Range("B10").Value2 = 'Afluria Quadrivalent'
$ws.Range("A11").Value2 = 'Influenza  (Age 5 years and older)'
$ws.Range("B11").Value2 = 'Afluria Quadrivalent'
